$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header text changes
$ws.Range("B1").Value = "Government-Cadw"
$ws.Range("C1").Value = "Government-Local_Authority"
$ws.Range("D1").Value = "Government-National"
$ws.Range("E1").Value = "Government-Other"
$ws.Range("F1").Value = "Independent-English_Heritage"
$ws.Range("G1").Value = "Independent-Historic_Environment_Scotland"
$ws.Range("H1").Value = "Independent-National_Trust"
$ws.Range("I1").Value = "Independent-National_Trust_for_Scotland"
$ws.Range("J1").Value = "Independent-Not_for_profit"
$ws.Range("K1").Value = "Independent-Private"
$ws.Range("L1").Value = "Independent-Unknown"
$ws.Range("M1").Value = "University"
$ws.Range("N1").Value = "Unknown"

# Row 2 (government) value changes
$ws.Range("B2").Value = 0.07099999999999999
$ws.Range("C2").Value = 21.802
$ws.Range("D2").Value = 1.939
$ws.Range("E2").Value = 0.236
$ws.Range("O2").Value = 24.048

# Row 3 (independent) value changes
$ws.Range("F3").Value = 1.253
$ws.Range("G3").Value = 0.497
$ws.Range("H3").Value = 4.375
$ws.Range("I3").Value = 0.638
$ws.Range("J3").Value = 41.003
$ws.Range("K3").Value = 17.758
$ws.Range("L3").Value = 5.226
$ws.Range("O3").Value = 70.75

# Row 4 (university) value changes
$ws.Range("M4").Value = 2.601
$ws.Range("O4").Value = 2.601

# Row 5 (unknown_gov) value changes
$ws.Range("N5").Value = 2.601
$ws.Range("O5").Value = 2.601

# Row 6 (COL_TOT) value changes
$ws.Range("B6").Value = 0.07099999999999999
$ws.Range("C6").Value = 21.802
$ws.Range("D6").Value = 1.939
$ws.Range("E6").Value = 0.236
$ws.Range("F6").Value = 1.253
$ws.Range("G6").Value = 0.497
$ws.Range("H6").Value = 4.375
$ws.Range("I6").Value = 0.638
$ws.Range("J6").Value = 41.003
$ws.Range("K6").Value = 17.758
$ws.Range("L6").Value = 5.226
$ws.Range("M6").Value = 2.601
$ws.Range("N6").Value = 2.601
